# 16.2.1.xlsx update
# - Updates a handful of header/label strings (adds "By"/"По" prefixes,
#   fixes capitalisation) on the "16.2.1" worksheet.
# - Adds placeholder "-" values for rows that previously had blank D cells
#   (functional-difficulties breakdown rows).
# - Wraps text and increases row height for the two sub-header rows that
#   now hold longer labels.
# - Leaves the active selection on B30, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-empty D cells with a "-" placeholder ---------
# (set first so the new shared string "-" lands before the relabelled
#  header cells below, matching the saved workbook's string order)
$ws.Range("D28").Value = "-"
$ws.Range("D29").Value = "-"
$ws.Range("D31").Value = "-"
$ws.Range("D32").Value = "-"
$ws.Range("D33").Value = "-"

# --- Relabel a few header/category cells -----------------------------
$ws.Range("C6").Value  = "By sex"
$ws.Range("C12").Value = "By territory"
$ws.Range("C22").Value = "By age (in years)"
$ws.Range("C34").Value = "Wealth quintile"
$ws.Range("B22").Value = "По возрасту (в годах)"
$ws.Range("A22").Value = "Жаш курагы боюнча (жылдарда)"

# --- Wrap text + taller rows for the two sub-header rows --------------
$ws.Range("A27:B27").WrapText = $true
$ws.Range("A30:B30").WrapText = $true
$ws.Rows.Item(27).RowHeight = 24
$ws.Rows.Item(30).RowHeight = 36

# --- Restore the saved selection/active cell --------------------------
$ws.Range("B30").Select()
